# Split the "Pay attention to their wants and needs" run on slide 2 into
# two runs: "Pay attention to their wants and " + "daily schedules",
# matching the commit "daily schedules for instructions".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shape = $s.Shapes.Item(2)   # "Content Placeholder 2"

$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 3 currently holds the single run:
#   "Pay attention to their wants and needs"
$prefix = "Pay attention to their wants and "
$oldTail = "needs"
$newTail = "daily schedules"

$para = $tr.Paragraphs(3, 1)
$run1 = $para.Runs(1, 1)

if ($run1.Text -ne ($prefix + $oldTail)) {
    throw "Unexpected source text for paragraph 3: [$($run1.Text)]"
}

# Trim the run's text back to the shared prefix, keeping its formatting.
$run1.Text = $prefix

# Re-acquire a fresh reference to that (now shorter) run and append a new
# run, inheriting the same character formatting, for the replacement tail.
$tr = $shape.TextFrame.TextRange
$para = $tr.Paragraphs(3, 1)
$run1 = $para.Runs(1, 1)
$run2 = $run1.InsertAfter($newTail)
